$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4705779254436493
$ws.Range("B1").Value = 1.342543244361877
$ws.Range("C1").Value = 6.550336837768555
$ws.Range("D1").Value = 2.002867460250854
$ws.Range("E1").Value = 1.745929837226868
